# add feature: truncate "reserved" fields, re-organize code structure
#
# REG1 (rows 4-7): fields 2 and 3 get collapsed/truncated into a single
# "reserved" field, and the stray port mapping in row 5 is cleared.
# REG2 (rows 9-13): the portB/portC port-bit mappings are re-organized.
# REG3 (rows 15-17): a brand new register block is appended after REG2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- REG1 field updates ----
$ws.Cells.Item(5, 2).Value = "reserved"      # B5: field2 -> reserved
$ws.Cells.Item(5, 6).ClearContents()         # F5: portB[ 4: 0] -> (blank)
$ws.Cells.Item(6, 5).Value = "RW"            # E6: RO -> RW

# ---- New REG3 register block (rows 15-17) ----
# Pull down the formatting used by the other register blocks: row 9 is the
# "register" header style, row 10 is the "field" style, row 14 is the blank
# trailing "comment" row style.
$ws.Range("A9:G9").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)
$ws.Range("A10:G10").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)
$ws.Range("A14:G14").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(15, 1).Value = "register"
$ws.Cells.Item(15, 2).Value = "REG3"
$ws.Cells.Item(15, 3).Value = "0x8"
$ws.Cells.Item(15, 7).Value = "reg desc"

$ws.Cells.Item(16, 1).Value = "field"
$ws.Cells.Item(16, 2).Value = "field1"
$ws.Cells.Item(16, 3).Value = "[ 31: 0]"
$ws.Cells.Item(16, 4).Value = "0x0"
$ws.Cells.Item(16, 5).Value = "RW"
$ws.Cells.Item(16, 6).Value = "portD[31:0]"
$ws.Cells.Item(16, 7).Value = "field desc"

$ws.Cells.Item(17, 1).Value = "comment"

$ws.Rows(15).RowHeight = 15.6
$ws.Rows(16).RowHeight = 15
$ws.Rows(17).RowHeight = 15

# ---- REG1 / REG2 field updates that introduce the last two new strings ----
$ws.Cells.Item(6, 6).Value = "portB[ 1: 0]"  # F6: portB[ 6: 5] -> portB[ 1: 0]
$ws.Cells.Item(11, 6).Value = "portB[2   ]"  # F11: portB[7   ] -> portB[2   ]

# ---- Selection / view bookkeeping ----
$ws.Range("F11").Select()
